$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10167.904
$ws.Range("I19").Value = 582.6667
$ws.Range("J19").Value = 22948.223
$ws.Range("K19").Value = 582.6667
$ws.Range("L19").Value = 22948.223
$ws.Range("M19").Value = -407.6667
$ws.Range("N19").Value = -23298.223
$ws.Range("H33").Value = 687.7222
$ws.Range("I33").Value = 507.37036
$ws.Range("J33").Value = 1228.7778
$ws.Range("K33").Value = 507.37036
$ws.Range("L33").Value = 1228.7778
$ws.Range("M33").Value = -278.37036
$ws.Range("N33").Value = -1686.7778
$ws.Range("H116").Value = 3237.2307
$ws.Range("I116").Value = 2763.8333
$ws.Range("J116").Value = 3643
$ws.Range("K116").Value = 2763.8333
$ws.Range("L116").Value = 3643
$ws.Range("M116").Value = 678.1667000000002
$ws.Range("N116").Value = -10527
$ws.Range("H129").Value = 1135.7222
$ws.Range("I129").Value = 558.1818
$ws.Range("J129").Value = 1283.4651
$ws.Range("K129").Value = 1674.5454
$ws.Range("L129").Value = 3850.3953
$ws.Range("M129").Value = 3325.4546
$ws.Range("N129").Value = -13850.3953
$ws.Range("H135").Value = 4309.5713
$ws.Range("I135").Value = 3194.5
$ws.Range("J135").Value = 11000
$ws.Range("K135").Value = 28750.5
$ws.Range("L135").Value = 99000
$ws.Range("M135").Value = -26215.5
$ws.Range("N135").Value = -104070
$ws.Range("H138").Value = 2266.2205
$ws.Range("J138").Value = 2128.6047
$ws.Range("L138").Value = 6385.8141
$ws.Range("N138").Value = -16665.8141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2262.1936
$ws.Range("I74").Value = 1862.5714
$ws.Range("J74").Value = 3101.4
$ws.Range("K74").Value = 1862.5714
$ws.Range("L74").Value = 3101.4
$ws.Range("M74").Value = -988.5714
$ws.Range("N74").Value = -4849.4
$ws.Range("H77").Value = 2262.1936
$ws.Range("I77").Value = 1862.5714
$ws.Range("J77").Value = 3101.4
$ws.Range("K77").Value = 9312.857
$ws.Range("L77").Value = 15507
$ws.Range("M77").Value = -4944.857
$ws.Range("N77").Value = -24243
$ws.Range("H97").Value = 1050
$ws.Range("I97").Value = 1050
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -554
$ws.Range("H132").Value = 2536.276
$ws.Range("I132").Value = 2140.2
$ws.Range("J132").Value = 3416.4443
$ws.Range("K132").Value = 6420.599999999999
$ws.Range("L132").Value = 10249.3329
$ws.Range("M132").Value = -3890.599999999999
$ws.Range("N132").Value = -15309.3329
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 627.1111
$ws.Range("I64").Value = 615.2
$ws.Range("J64").Value = 642
$ws.Range("K64").Value = 615.2
$ws.Range("L64").Value = 642
$ws.Range("M64").Value = -390.2
$ws.Range("N64").Value = -1092
$ws.Range("H67").Value = 627.1111
$ws.Range("I67").Value = 615.2
$ws.Range("J67").Value = 642
$ws.Range("K67").Value = 615.2
$ws.Range("L67").Value = 642
$ws.Range("M67").Value = 164.8
$ws.Range("N67").Value = -2202
$ws.Range("H81").Value = 64207.145
$ws.Range("J81").Value = 64207.145
$ws.Range("L81").Value = 64207.145
$ws.Range("N81").Value = -66329.14499999999
$ws.Range("H84").Value = 64207.145
$ws.Range("J84").Value = 64207.145
$ws.Range("L84").Value = 192621.435
$ws.Range("N84").Value = -203229.435
$ws.Range("H86").Value = 52633676
$ws.Range("I86").Value = 58825612
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 58825612
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -58824489
$ws.Range("N86").Value = -4496
$ws.Range("H89").Value = 52633676
$ws.Range("I89").Value = 58825612
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 294128060
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -294122444
$ws.Range("N89").Value = -22482
$ws.Range("H95").Value = 55000.5
$ws.Range("J95").Value = 55000.5
$ws.Range("L95").Value = 55000.5
$ws.Range("N95").Value = -60492.5
$ws.Range("H96").Value = 10685.4
$ws.Range("I96").Value = 10685.4
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 10685.4
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = -7939.4
$ws.Range("H97").Value = 29998.334
$ws.Range("J97").Value = 34997.5
$ws.Range("L97").Value = 34997.5
$ws.Range("N97").Value = -36979.5
$ws.Range("H99").Value = 1790
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1790
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = 1790
$ws.Range("N99").Value = -4786
$ws.Range("H135").Value = 39246.668
$ws.Range("J135").Value = 39246.668
$ws.Range("L135").Value = 39246.668
$ws.Range("N135").Value = -49386.668
$ws.Range("N96").ClearContents()
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2268.3
$ws.Range("I16").Value = 1788.8
$ws.Range("J16").Value = 2747.8
$ws.Range("K16").Value = 1788.8
$ws.Range("L16").Value = 2747.8
$ws.Range("M16").Value = -1501.8
$ws.Range("N16").Value = -3321.8
$ws.Range("H31").Value = 6678.8184
$ws.Range("I31").Value = 1672.5714
$ws.Range("J31").Value = 11249.739
$ws.Range("K31").Value = 1672.5714
$ws.Range("L31").Value = 11249.739
$ws.Range("M31").Value = -1377.5714
$ws.Range("N31").Value = -11839.739
$ws.Range("H34").Value = 6678.8184
$ws.Range("I34").Value = 1672.5714
$ws.Range("J34").Value = 11249.739
$ws.Range("K34").Value = 1672.5714
$ws.Range("L34").Value = 11249.739
$ws.Range("M34").Value = -1470.5714
$ws.Range("N34").Value = -11653.739
$ws.Range("H113").Value = 2268.3
$ws.Range("I113").Value = 1788.8
$ws.Range("J113").Value = 2747.8
$ws.Range("K113").Value = 1788.8
$ws.Range("L113").Value = 2747.8
$ws.Range("M113").Value = 381.2
$ws.Range("N113").Value = -7087.8
$ws.Range("H134").Value = 2486.7144
$ws.Range("I134").Value = 2139.25
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 6417.75
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -3882.75
$ws.Range("N134").Value = -13920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2766.6667
$ws.Range("J35").Value = 4000
$ws.Range("L35").Value = 12000
$ws.Range("N35").Value = -12576
$ws.Range("H39").Value = 1566.2461
$ws.Range("I39").Value = 100
$ws.Range("J39").Value = 1589.1562
$ws.Range("K39").Value = 300
$ws.Range("L39").Value = 4767.4686
$ws.Range("M39").Value = -6
$ws.Range("N39").Value = -5355.4686
$ws.Range("H68").Value = 781.8182
$ws.Range("I68").Value = 757.1429000000001
$ws.Range("J68").Value = 825
$ws.Range("K68").Value = 2271.4287
$ws.Range("L68").Value = 2475
$ws.Range("M68").Value = -1460.4287
$ws.Range("N68").Value = -4097
$ws.Range("H71").Value = 781.8182
$ws.Range("I71").Value = 757.1429000000001
$ws.Range("J71").Value = 825
$ws.Range("K71").Value = 6814.2861
$ws.Range("L71").Value = 7425
$ws.Range("M71").Value = -2758.2861
$ws.Range("N71").Value = -15537
$ws.Range("H80").Value = 4497.6665
$ws.Range("J80").Value = 4497.6665
$ws.Range("L80").Value = 13492.9995
$ws.Range("N80").Value = -15364.9995
$ws.Range("H83").Value = 4497.6665
$ws.Range("J83").Value = 4497.6665
$ws.Range("L83").Value = 40478.9985
$ws.Range("N83").Value = -49838.9985
$ws.Range("H131").Value = 885.55554
$ws.Range("J131").Value = 1060
$ws.Range("L131").Value = 3180
$ws.Range("N131").Value = -13260
$ws.Range("H134").Value = 4817.2905
$ws.Range("I134").Value = 2464.6428
$ws.Range("J134").Value = 6754.7646
$ws.Range("K134").Value = 7393.928400000001
$ws.Range("L134").Value = 20264.2938
$ws.Range("M134").Value = -2323.928400000001
$ws.Range("N134").Value = -30404.2938
$ws.Range("H141").Value = 5618.793
$ws.Range("I141").Value = 4338.75
$ws.Range("J141").Value = 7194.231
$ws.Range("K141").Value = 13016.25
$ws.Range("L141").Value = 21582.693
$ws.Range("M141").Value = -7836.25
$ws.Range("N141").Value = -31942.693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2297.5483
$ws.Range("I68").Value = 1901.8462
$ws.Range("J68").Value = 2583.3333
$ws.Range("K68").Value = 1901.8462
$ws.Range("L68").Value = 2583.3333
$ws.Range("M68").Value = -1152.8462
$ws.Range("N68").Value = -4081.3333
$ws.Range("H71").Value = 2297.5483
$ws.Range("I71").Value = 1901.8462
$ws.Range("J71").Value = 2583.3333
$ws.Range("K71").Value = 9509.231
$ws.Range("L71").Value = 12916.6665
$ws.Range("M71").Value = -5765.231
$ws.Range("N71").Value = -20404.6665
$ws.Range("H95").Value = 90000
$ws.Range("J95").Value = 90000
$ws.Range("L95").Value = 90000
$ws.Range("N95").Value = -95492
$ws.Range("H136").Value = 7248347.5
$ws.Range("I136").Value = 1526.625
$ws.Range("J136").Value = 23812510
$ws.Range("K136").Value = 4579.875
$ws.Range("L136").Value = 71437530
$ws.Range("M136").Value = -2029.875
$ws.Range("N136").Value = -71442630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 69124.5
$ws.Range("J63").Value = 69124.5
$ws.Range("L63").Value = 69124.5
$ws.Range("N63").Value = -70372.5
$ws.Range("H66").Value = 69124.5
$ws.Range("J66").Value = 69124.5
$ws.Range("L66").Value = 207373.5
$ws.Range("N66").Value = -213613.5
$ws.Range("H96").Value = 4458.7334
$ws.Range("I96").Value = 3478.2
$ws.Range("J96").Value = 6419.8
$ws.Range("K96").Value = 3478.2
$ws.Range("L96").Value = 6419.8
$ws.Range("M96").Value = -2105.2
$ws.Range("N96").Value = -9165.799999999999
$ws.Range("H132").Value = 4275465
$ws.Range("I132").Value = 1850.8334
$ws.Range("J132").Value = 6174849
$ws.Range("K132").Value = 5552.5002
$ws.Range("L132").Value = 18524547
$ws.Range("M132").Value = -3022.5002
$ws.Range("N132").Value = -18529607
